# thêm chức năng quay lại giao diện chính (lựa chọn nhân sự / phòng ban)
#
# Data-level effect of this commit on the embedded "database" workbook:
#   - The "PB02" department row is removed from the PhongBan sheet and the
#     rows below it shift up (dimension A1:F6 -> A1:F5).
#   - The employees that used to be linked to PB02 are reassigned to other
#     departments: one to PB01, two to PB05 - so PB01's headcount (2 -> 3)
#     and PB05's headcount (2 -> 4) grow while the overall total headcount
#     stays the same (11).

$wb = $excel.ActiveWorkbook
$wsPhongBan = $wb.Worksheets.Item("PhongBan")
$wsNhanSu   = $wb.Worksheets.Item("NhanSu")

# --- PhongBan: remove the PB02 row (row 3) entirely; rows 4-6 shift up to
#     become rows 3-5.
$wsPhongBan.Rows(3).Delete()

# Update tongSoNhanVien (headcount) for the departments that absorbed the
# reassigned employees.
$wsPhongBan.Range("F2").Value = 3   # PB01: 2 -> 3
$wsPhongBan.Range("F5").Value = 4   # PB05: 2 -> 4 (was row 6 before the delete)

# --- NhanSu: re-point the employees that used to belong to PB02 at their
#     new department.
$wsNhanSu.Range("H3").Value  = "PB01"
$wsNhanSu.Range("H7").Value  = "PB05"
$wsNhanSu.Range("H12").Value = "PB05"
